# Adds the new "SHA3-SHAKE Registration" sheet, populates it with the
# algorithm property table (mirroring the existing "Symmetric Registration"
# / "Symmetric Prompt Test Group" sheets), and updates the selections /
# active-tab state on the existing sheets to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the selection on "Symmetric Registration" (sheet 1) first,
#    while it is still the active tab.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1:G2").Select()

# ---------------------------------------------------------------------
# 2. Append a brand-new worksheet after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "SHA3-SHAKE Registration"

# Column widths (character units); ColumnWidth assignments round-trip
# through the host's internal MDW-based storage, so these are chosen to
# land close to the original bestFit widths.
$ws3.Columns.Item(1).ColumnWidth = 12.1666666666667
$ws3.Columns.Item(2).ColumnWidth = 9.73333333333333
$ws3.Columns.Item(3).ColumnWidth = 9.73333333333333
$ws3.Columns.Item(4).ColumnWidth = 34.1666666666667
$ws3.Columns.Item(5).ColumnWidth = 9.73333333333333
$ws3.Columns.Item(6).ColumnWidth = 17.8833333333333

# Row 1 - header row
$ws3.Range("A1").Value = "algorithm"
$ws3.Range("B1").Value = "inBit"
$ws3.Range("C1").Value = "inEmpty"
$ws3.Range("D1").Value = "outLength"
$ws3.Range("E1").Value = "outBit"
$ws3.Range("F1").Value = '<ttcol align="left">'
$ws3.Range("G1").Value = "</ttcol>"
$ws3.Range("H1").Formula = "=CONCAT(F1,TEXTJOIN(CONCAT(G1,F1),,A1:E1),G1)"

# Row 2 - SHA3-224
$ws3.Range("A2").Value = "SHA3-224"
$ws3.Range("B2").Value = "true, false"
$ws3.Range("C2").Value = "true, false"
$ws3.Range("F2").Value = "<c>"
$ws3.Range("G2").Value = "</c>"
$ws3.Range("H2").Formula = "=CONCAT(F2,TEXTJOIN(CONCAT(G2,F2),FALSE,A2:E2),G2)"

# Row 3 - SHA3-256
$ws3.Range("A3").Value = "SHA3-256"
$ws3.Range("B3").Value = "true, false"
$ws3.Range("C3").Value = "true, false"
$ws3.Range("F3").Value = "<c>"
$ws3.Range("G3").Value = "</c>"
$ws3.Range("H3").Formula = "=CONCAT(F3,TEXTJOIN(CONCAT(G3,F3),FALSE,A3:E3),G3)"

# Row 4 - SHA3-384
$ws3.Range("A4").Value = "SHA3-384"
$ws3.Range("B4").Value = "true, false"
$ws3.Range("C4").Value = "true, false"
$ws3.Range("F4").Value = "<c>"
$ws3.Range("G4").Value = "</c>"
$ws3.Range("H4").Formula = "=CONCAT(F4,TEXTJOIN(CONCAT(G4,F4),FALSE,A4:E4),G4)"

# Row 5 - SHA3-512
$ws3.Range("A5").Value = "SHA3-512"
$ws3.Range("B5").Value = "true, false"
$ws3.Range("C5").Value = "true, false"
$ws3.Range("F5").Value = "<c>"
$ws3.Range("G5").Value = "</c>"
$ws3.Range("H5").Formula = "=CONCAT(F5,TEXTJOIN(CONCAT(G5,F5),FALSE,A5:E5),G5)"

# Row 6 - SHAKE-128
$ws3.Range("A6").Value = "SHAKE-128"
$ws3.Range("B6").Value = "true, false"
$ws3.Range("C6").Value = "true, false"
$ws3.Range("D6").Value = '{"Min": 16, "Max": 65536, "Inc": any}'
$ws3.Range("E6").Value = "true, false"
$ws3.Range("F6").Value = "<c>"
$ws3.Range("G6").Value = "</c>"
$ws3.Range("H6").Formula = "=CONCAT(F6,TEXTJOIN(CONCAT(G6,F6),FALSE,A6:E6),G6)"

# Row 7 - SHAKE-256
$ws3.Range("A7").Value = "SHAKE-256"
$ws3.Range("B7").Value = "true, false"
$ws3.Range("C7").Value = "true, false"
$ws3.Range("D7").Value = '{"Min": 16, "Max": 65536, "Inc": any}'
$ws3.Range("E7").Value = "true, false"
$ws3.Range("F7").Value = "<c>"
$ws3.Range("G7").Value = "</c>"
$ws3.Range("H7").Formula = "=CONCAT(F7,TEXTJOIN(CONCAT(G7,F7),FALSE,A7:E7),G7)"

# ---------------------------------------------------------------------
# 3. Put the selection on the new sheet where the author left it, and
#    make it the active tab (last Select() call wins / matches the
#    xr:uid-less default behaviour used throughout this workbook).
# ---------------------------------------------------------------------
$ws3.Range("I30").Select()
